# Generate Report for Handback
#
# The first localization file (6e664e64-...) has come back from handback
# and is in sync with en-US again, so its status flips from "Ready for
# handoff" to "Handed back: in sync with en-US" everywhere that status is
# shown (Overview + the per-language detail sheets). The per-language
# sheets also now carry the handback target/result files ("Latest Target
# File" / "Latest Handback File" columns) with their hyperlinks, and the
# handback timestamp ("Latest Handback DateTime").

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e36f3c5820925b29bf00e1cc1c0922a07fd1bbf1/e2e/6e664e64-057f-4d21-815c-ade30daf2d25.md",
    "",
    "",
    "6e664e64-057f-4d21-815c-ade30daf2d25.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ed0199ebb4d072d387d271307ad0753274597ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/6e664e64-057f-4d21-815c-ade30daf2d25.ae29106d17ee0f15f872328c333da460d432b114.zh-cn.xlf",
    "",
    "",
    "6e664e64-057f-4d21-815c-ade30daf2d25.ae29106d17ee0f15f872328c333da460d432b114.zh-cn.xlf"
) | Out-Null

$wsZh.Range("H2").Value = "2016-03-21 08:14:13"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e36f3c5820925b29bf00e1cc1c0922a07fd1bbf1/e2e/c876207b-9d86-4547-8eb5-bf7f68b50fd1.md",
    "",
    "",
    "c876207b-9d86-4547-8eb5-bf7f68b50fd1.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ed0199ebb4d072d387d271307ad0753274597ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c876207b-9d86-4547-8eb5-bf7f68b50fd1.ef22e670e00e246da637c573d5facb74300b3cfc.zh-cn.xlf",
    "",
    "",
    "c876207b-9d86-4547-8eb5-bf7f68b50fd1.ef22e670e00e246da637c573d5facb74300b3cfc.zh-cn.xlf"
) | Out-Null

$wsZh.Range("H3").Value = "2016-03-21 08:14:13"

# ---- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e36f3c5820925b29bf00e1cc1c0922a07fd1bbf1/e2e/6e664e64-057f-4d21-815c-ade30daf2d25.md",
    "",
    "",
    "6e664e64-057f-4d21-815c-ade30daf2d25.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67702089ccdb9a1cb60547f427fbc45d32f8e016/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/6e664e64-057f-4d21-815c-ade30daf2d25.ae29106d17ee0f15f872328c333da460d432b114.de-de.xlf",
    "",
    "",
    "6e664e64-057f-4d21-815c-ade30daf2d25.ae29106d17ee0f15f872328c333da460d432b114.de-de.xlf"
) | Out-Null

$wsDe.Range("H2").Value = "2016-03-21 08:14:18"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e36f3c5820925b29bf00e1cc1c0922a07fd1bbf1/e2e/c876207b-9d86-4547-8eb5-bf7f68b50fd1.md",
    "",
    "",
    "c876207b-9d86-4547-8eb5-bf7f68b50fd1.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67702089ccdb9a1cb60547f427fbc45d32f8e016/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c876207b-9d86-4547-8eb5-bf7f68b50fd1.ef22e670e00e246da637c573d5facb74300b3cfc.de-de.xlf",
    "",
    "",
    "c876207b-9d86-4547-8eb5-bf7f68b50fd1.ef22e670e00e246da637c573d5facb74300b3cfc.de-de.xlf"
) | Out-Null

$wsDe.Range("H3").Value = "2016-03-21 08:14:18"
